# Append 3 new captured rows to the "Captured_Values" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Captured_Values")

$text = "Real Programmers Count 0123456789 From Zero"

for ($i = 15; $i -le 17; $i++) {
    $ws.Cells.Item($i, 1).Value = 123456789
    $ws.Cells.Item($i, 2).Value = $text
}
